$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (below the current row 2); it inherits the
# formatting of row 2 above it, so the ID column stays unstyled and the date
# columns keep style index 2 (the date number format), matching the table's
# existing rows.
$ws.Rows.Item(3).Insert()

# Move the old row 2 (week 44) contents down into the newly inserted row 3.
$ws.Cells.Item(3, 1).Value = $ws.Cells.Item(2, 1).Value2
$ws.Cells.Item(3, 2).Value = $ws.Cells.Item(2, 2).Value2
$ws.Cells.Item(3, 3).Value = $ws.Cells.Item(2, 3).Value2

# Write the new week (45) into row 2, continuing the existing sequence.
$ws.Cells.Item(2, 1).Value = 45
$ws.Cells.Item(2, 2).Value = 45380
$ws.Cells.Item(2, 3).Value = 45409
